$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns are treated as text so values
# like "1.001" or "1.000" are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '22.053.44'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '1.557.69'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '1.000'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('D6').Value = '291.77'
$ws.Range('E6').Value = '  +1.64%  '
$ws.Range('D7').Value = '0.3966'
$ws.Range('E7').Value = '  +3.67%  '
$ws.Range('D8').Value = '0.3242'
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('D9').Value = '44.33'
$ws.Range('E9').Value = '  +1.74%  '
$ws.Range('D10').Value = '0.07273'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('D11').Value = '1.081'
$ws.Range('E11').Value = '  -3.91%  '
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = '5.715'
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('D14').Value = '18.85'
$ws.Range('E14').Value = '  -6.09%  '
$ws.Range('D15').Value = '6.663'
$ws.Range('E15').Value = '  -1.27%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.00001123'
$ws.Range('E16').Value = '  +3.83%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '1.554.96'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').Value = '0.06590'
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').Value = '83.88'
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('D20').Value = '0.9993'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').Value = '6.281'
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('D22').Value = '15.60'
$ws.Range('E22').Value = '  -2.65%  '
$ws.Range('D23').Value = '11.35'
$ws.Range('E23').Value = '  -2.78%  '
$ws.Range('D24').Value = '22.066.93'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = '2.365'
$ws.Range('E25').Value = '  +3.10%  '
$ws.Range('D26').Value = '2.428'
$ws.Range('E26').Value = '  -2.84%  '
$ws.Range('D27').Value = '148.43'
$ws.Range('D28').Value = '18.63'
$ws.Range('E28').Value = '  -2.63%  '
$ws.Range('D29').Value = '4.874'
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('D30').Value = '1.739.29'
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('D31').Value = '119.42'
$ws.Range('E31').Value = '  -1.79%  '
$ws.Range('D32').Value = '0.9897'
$ws.Range('E32').Value = '  -8.21%  '
$ws.Range('D33').Value = '5.942'
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('D34').Value = '0.08336'
$ws.Range('E34').Value = '  +1.52%  '
$ws.Range('D35').Value = '9.179'
$ws.Range('D36').Value = '1.613'
$ws.Range('E36').Value = '  -15.07%  '
$ws.Range('D37').Value = '0.02282'
$ws.Range('E37').Value = '  -1.29%  '
$ws.Range('D38').Value = '5.152'
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('D39').Value = '0.06024'
$ws.Range('E39').Value = '  -4.20%  '
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').Value = '0.2061'
$ws.Range('E41').Value = '  -3.88%  '
$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D42').Value = '0.9997'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '10.81'
$ws.Range('E43').Value = '  -1.85%  '
$ws.Range('D44').Value = '0.5845'
$ws.Range('E44').Value = '  -2.64%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '13.15'
$ws.Range('E45').Value = '  -3.55%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.768'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').Value = '0.5608'
$ws.Range('E47').Value = '  -3.55%  '
$ws.Range('D48').Value = '118.41'
$ws.Range('E48').Value = '  -2.67%  '
$ws.Range('D49').Value = '1.902'
$ws.Range('E49').Value = '  -3.26%  '
$ws.Range('D50').Value = '1.142'
$ws.Range('E50').Value = '  -2.50%  '
$ws.Range('D51').Value = '0.06822'
$ws.Range('E51').Value = '  -2.70%  '
